# Update "想去人数" (number of people interested) values in column F
# across the 4 worksheets, as per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1319
$ws.Range("F3").Value = 75
$ws.Range("F5").Value = 162
$ws.Range("F6").Value = 382
$ws.Range("F7").Value = 170
$ws.Range("F9").Value = 977
$ws.Range("F10").Value = 312
$ws.Range("F11").Value = 170
$ws.Range("F15").Value = 339
$ws.Range("F16").Value = 748
$ws.Range("F17").Value = 124
$ws.Range("F18").Value = 701
$ws.Range("F19").Value = 245
$ws.Range("F20").Value = 62
$ws.Range("F21").Value = 963
$ws.Range("F22").Value = 424
$ws.Range("F23").Value = 234
$ws.Range("F24").Value = 75
$ws.Range("F25").Value = 347
$ws.Range("F28").Value = 448

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 352
$ws.Range("F11").Value = 143

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 341

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 341
$ws.Range("F3").Value = 1319
$ws.Range("F4").Value = 75
$ws.Range("F7").Value = 162
$ws.Range("F8").Value = 382
$ws.Range("F9").Value = 170
$ws.Range("F11").Value = 977
$ws.Range("F12").Value = 312
$ws.Range("F13").Value = 170
$ws.Range("F16").Value = 352
$ws.Range("F22").Value = 339
$ws.Range("F23").Value = 748
$ws.Range("F24").Value = 124
$ws.Range("F25").Value = 701
$ws.Range("F26").Value = 245
$ws.Range("F27").Value = 62
$ws.Range("F28").Value = 963
$ws.Range("F29").Value = 424
$ws.Range("F32").Value = 234
$ws.Range("F33").Value = 75
$ws.Range("F34").Value = 347
$ws.Range("F36").Value = 143
$ws.Range("F40").Value = 448
